$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '45.044.53'
$ws.Range('E2').Value = '  +4.42%  '
$ws.Range('D3').Value = '2.426.93'
$ws.Range('E3').Value = '  +2.31%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '318.25'
$ws.Range('E5').Value = '  +4.72%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '104.31'
$ws.Range('E6').Value = '  +9.18%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.516'
$ws.Range('E7').Value = '  +2.62%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.530'
$ws.Range('E9').Value = '  +10.19%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.73'
$ws.Range('E10').Value = '  +3.79%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0803'
$ws.Range('E11').Value = '  +2.01%  '
$ws.Range('E12').Value = '  -2.59%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.55'
$ws.Range('E13').Value = '  +1.89%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.96'
$ws.Range('E14').Value = '  +2.61%  '
$ws.Range('D15').Value = '2.810.99'
$ws.Range('E15').Value = '  +2.63%  '
$ws.Range('D16').Value = '2.428.54'
$ws.Range('E16').Value = '  +2.79%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.835'
$ws.Range('E17').Value = '  +4.46%  '
$ws.Range('D18').Value = '44.925.05'
$ws.Range('E18').Value = '  +4.19%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.39'
$ws.Range('E19').Value = '  +3.48%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.35'
$ws.Range('E20').Value = '  +1.45%  '
$ws.Range('E21').Value = '  +3.64%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.83'
$ws.Range('E22').Value = '  +1.31%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '243.58'
$ws.Range('E23').Value = '  +3.54%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.29'
$ws.Range('E24').Value = '  +4.13%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.51'
$ws.Range('E25').Value = '  +2.69%  '
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.39'
$ws.Range('E27').Value = '  +3.63%  '
$ws.Range('E28').Value = '  -7.55%  '
$ws.Range('E29').Value = '  +2.50%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '33.89'
$ws.Range('E30').Value = '  +5.48%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '48.89'
$ws.Range('E31').Value = '  +1.85%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.129'
$ws.Range('E32').Value = '  +17.31%  '
$ws.Range('E33').Value = '  +12.08%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.22'
$ws.Range('E34').Value = '  +3.90%  '
$ws.Range('E35').Value = '  +0.28%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0764'
$ws.Range('E36').Value = '  +4.04%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.90'
$ws.Range('E37').Value = '  +4.29%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.51'
$ws.Range('E38').Value = '  +4.53%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '127.10'
$ws.Range('E39').Value = '  -1.47%  '
$ws.Range('E40').Value = '  +0.94%  '
$ws.Range('E41').Value = '  +2.16%  '
$ws.Range('E42').Value = '  -2.83%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '21.10'
$ws.Range('E43').Value = '  -0.16%  '
$ws.Range('E44').Value = '  +4.44%  '
$ws.Range('D45').Value = '1.943.32'
$ws.Range('E45').Value = '  +0.67%  '
$ws.Range('E46').Value = '  -0.53%  '
$ws.Range('E47').Value = '  +8.28%  '
$ws.Range('E48').Value = '  -0.24%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.78'
$ws.Range('E49').Value = '  +17.91%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '75.80'
$ws.Range('E50').Value = '  +6.13%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '54.11'
$ws.Range('E51').Value = '  +5.13%  '
